$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.963.32'
$ws.Range('E2').Value = '  +2.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.045.15'
$ws.Range('E3').Value = '  +2.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.54'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.84'
$ws.Range('E6').Value = '  +7.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.040.68'
$ws.Range('E8').Value = '  +2.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.42'
$ws.Range('E10').Value = '  +11.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').Value = '  +6.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('E12').Value = '  +2.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000235'
$ws.Range('E13').Value = '  +4.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.05'
$ws.Range('E14').Value = '  +3.55%  '
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.539.06'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.930.74'
$ws.Range('E17').Value = '  +2.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.05'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.031.88'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.45'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.28'
$ws.Range('E21').Value = '  +3.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.695'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.51'
$ws.Range('E23').Value = '  +2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.65'
$ws.Range('E24').Value = '  +1.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.29'
$ws.Range('E25').Value = '  +6.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.94'
$ws.Range('E26').Value = '  +10.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.27'
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.73'
$ws.Range('E29').Value = '  +3.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.44'
$ws.Range('E30').Value = '  +9.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.17'
$ws.Range('E32').Value = '  +5.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.72'
$ws.Range('E33').Value = '  +3.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.111'
$ws.Range('E34').Value = '  +5.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0868'
$ws.Range('E35').Value = '  +12.00%  '
$ws.Range('E36').Value = '  +2.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.90'
$ws.Range('E37').Value = '  +3.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.16'
$ws.Range('E38').Value = '  +15.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.11'
$ws.Range('E39').Value = '  +3.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.44'
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.10'
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('E42').Value = '  +4.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.295'
$ws.Range('E43').Value = '  +12.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.53'
$ws.Range('E44').Value = '  +12.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '394.58'
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0357'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.750.25'
$ws.Range('E47').Value = '  +2.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.41'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.21'
$ws.Range('E50').Value = '  +3.38%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.17'
$ws.Range('E51').Value = '  +3.98%  '
